# Finish the "Packing Slip" tutorial worksheet: add a sales-tax rate,
# a computed tax line, a shipping line, and a grand total.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packing Slip")

# --- Row 9: sales tax rate -------------------------------------------------
$ws.Range("F9").Value = "Sales Tax"
$ws.Range("G9").Value = 0.05
$ws.Range("G9").NumberFormat = "0%"

# Row 10 (item header row) grows taller to fit the extra info above it.
$ws.Rows.Item(10).RowHeight = 27

# --- Row 17: tax amount, computed from subtotal (G16) * tax rate (G9) -----
$ws.Range("F17").Value = "Tax"
$ws.Range("G17").Formula = "=G16*G9"
$ws.Range("G17").NumberFormat = "`"$`"#,##0.00"

# --- Row 18: shipping charge ------------------------------------------------
$ws.Range("F18").Value = "Shipping"
$ws.Range("G18").Value = 45
$ws.Range("G18").NumberFormat = "`"$`"#,##0.00_);[Red]\(`"$`"#,##0.00\)"

# --- Row 19: grand total -----------------------------------------------------
$ws.Range("F19").Value = "Total"
$ws.Range("G19").Formula = "=SUM(G16:G18)"
$ws.Range("G19").NumberFormat = "`"$`"#,##0.00"

# The selection marker left over from editing is no longer meaningful;
# drop it so the sheet view just shows the default top-left state.
$ws.Range("A1").Select()
